$d = $word.ActiveDocument

# Replace 1: Predicate (Role) line (appears twice)
$d.Content.Find.Execute(
    "Predicate (Role): AggregationSubjectKind. Employee / Employer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Predicate (Role): AggregationSubjectKind. Employee / Employer / Position",
    2)

# Replace 2: Object (Occurrence) line
$d.Content.Find.Execute(
    "Object (Occurrence): AggregationSubject. anEmployee / anEmployer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Object (Occurrence): AggregationSubject. anEmployee / anEmployer / aPosition",
    2)

# Replace 3: Subject (Context) line
$d.Content.Find.Execute(
    "Subject (Context): AggregationSubject. anEmployee / anEmployer",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Subject (Context): AggregationSubject. anEmployee / anEmployer / aPosition",
    2)
